$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The R11 (essround 11) block of country rows is missing four countries
# (BG, IL, LV, ME). Insert a blank row at each alphabetically-correct
# position (processed top-to-bottom so each target row index already
# accounts for the earlier inserts) and fill in the new row's values.

$ws.Rows(109).Insert()
$ws.Range("A109").Value = "R11"
$ws.Range("B109").Value = "BG"
$ws.Range("C109").Value = 2239
$ws.Range("D109").Value = 2238.99997510016
$ws.Range("E109").Value = 2239.00000280142
$ws.Range("F109").Value = 5534290.98106921
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 0

$ws.Rows(122).Insert()
$ws.Range("A122").Value = "R11"
$ws.Range("B122").Value = "IL"
$ws.Range("C122").Value = 906
$ws.Range("D122").Value = 905.999990701675
$ws.Range("E122").Value = 905.999990582466
$ws.Range("F122").Value = 6750299.96842146
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 0

$ws.Rows(126).Insert()
$ws.Range("A126").Value = "R11"
$ws.Range("B126").Value = "LV"
$ws.Range("C126").Value = 1252
$ws.Range("D126").Value = 1252.00000301003
$ws.Range("E126").Value = 1252.00000101328
$ws.Range("F126").Value = 1582531.00175411
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 0

$ws.Rows(127).Insert()
$ws.Range("A127").Value = "R11"
$ws.Range("B127").Value = "ME"
$ws.Range("C127").Value = 1609
$ws.Range("D127").Value = 1609.0000321418
$ws.Range("E127").Value = 1609.00002001971
$ws.Range("F127").Value = 506465.000677854
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 0
